$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells for the 5 per-tag errors + the overall average ---
$ws.Range("V1").Value = "Error1"
$ws.Range("W1").Value = "Error2"
$ws.Range("X1").Value = "Error3"
$ws.Range("Y1").Value = "Error4"
$ws.Range("Z1").Value = "Error5"
$ws.Range("AA1").Value = "ErrorAverage"

# --- Column widths for the new columns (matches narrower Error columns + wider ErrorAverage column) ---
$ws.Range("V1:Z1").ColumnWidth = 5.666666666666667
$ws.Range("AA1").ColumnWidth = 11.833333333333332

# --- Refresh the randomized tag coordinates (columns A:F) for every data row ---
    $ws.Cells.Item(2, 1).Value = 78
    $ws.Cells.Item(2, 2).Value = 195
    $ws.Cells.Item(2, 3).Value = 58
    $ws.Cells.Item(2, 4).Value = 117
    $ws.Cells.Item(2, 5).Value = 112
    $ws.Cells.Item(2, 6).Value = 121
    $ws.Cells.Item(3, 1).Value = 113
    $ws.Cells.Item(3, 2).Value = 56
    $ws.Cells.Item(3, 3).Value = 31
    $ws.Cells.Item(3, 4).Value = 145
    $ws.Cells.Item(3, 5).Value = 8
    $ws.Cells.Item(3, 6).Value = 157
    $ws.Cells.Item(4, 1).Value = 82
    $ws.Cells.Item(4, 2).Value = 166
    $ws.Cells.Item(4, 3).Value = 193
    $ws.Cells.Item(4, 4).Value = 84
    $ws.Cells.Item(4, 5).Value = 11
    $ws.Cells.Item(4, 6).Value = 154
    $ws.Cells.Item(5, 1).Value = 127
    $ws.Cells.Item(5, 2).Value = 102
    $ws.Cells.Item(5, 3).Value = 113
    $ws.Cells.Item(5, 4).Value = 61
    $ws.Cells.Item(5, 5).Value = 147
    $ws.Cells.Item(5, 6).Value = 66
    $ws.Cells.Item(6, 1).Value = 177
    $ws.Cells.Item(6, 2).Value = 67
    $ws.Cells.Item(6, 3).Value = 158
    $ws.Cells.Item(6, 4).Value = 159
    $ws.Cells.Item(6, 5).Value = 10
    $ws.Cells.Item(6, 6).Value = 79
    $ws.Cells.Item(7, 1).Value = 194
    $ws.Cells.Item(7, 2).Value = 83
    $ws.Cells.Item(7, 3).Value = 66
    $ws.Cells.Item(7, 4).Value = 126
    $ws.Cells.Item(7, 5).Value = 169
    $ws.Cells.Item(7, 6).Value = 154
    $ws.Cells.Item(8, 1).Value = 108
    $ws.Cells.Item(8, 2).Value = 148
    $ws.Cells.Item(8, 3).Value = 147
    $ws.Cells.Item(8, 4).Value = 88
    $ws.Cells.Item(8, 5).Value = 167
    $ws.Cells.Item(8, 6).Value = 97
    $ws.Cells.Item(9, 1).Value = 190
    $ws.Cells.Item(9, 2).Value = 159
    $ws.Cells.Item(9, 3).Value = 27
    $ws.Cells.Item(9, 4).Value = 47
    $ws.Cells.Item(9, 5).Value = 26
    $ws.Cells.Item(9, 6).Value = 105
    $ws.Cells.Item(10, 1).Value = 12
    $ws.Cells.Item(10, 2).Value = 94
    $ws.Cells.Item(10, 3).Value = 195
    $ws.Cells.Item(10, 4).Value = 135
    $ws.Cells.Item(10, 5).Value = 171
    $ws.Cells.Item(10, 6).Value = 97
    $ws.Cells.Item(11, 1).Value = 58
    $ws.Cells.Item(11, 2).Value = 37
    $ws.Cells.Item(11, 3).Value = 191
    $ws.Cells.Item(11, 4).Value = 149
    $ws.Cells.Item(11, 5).Value = 187
    $ws.Cells.Item(11, 6).Value = 131
    $ws.Cells.Item(12, 1).Value = 38
    $ws.Cells.Item(12, 2).Value = 58
    $ws.Cells.Item(12, 3).Value = 12
    $ws.Cells.Item(12, 4).Value = 9
    $ws.Cells.Item(12, 5).Value = 55
    $ws.Cells.Item(12, 6).Value = 34
    $ws.Cells.Item(13, 1).Value = 136
    $ws.Cells.Item(13, 2).Value = 180
    $ws.Cells.Item(13, 3).Value = 70
    $ws.Cells.Item(13, 4).Value = 20
    $ws.Cells.Item(13, 5).Value = 156
    $ws.Cells.Item(13, 6).Value = 92
    $ws.Cells.Item(14, 1).Value = 138
    $ws.Cells.Item(14, 2).Value = 136
    $ws.Cells.Item(14, 3).Value = 84
    $ws.Cells.Item(14, 4).Value = 130
    $ws.Cells.Item(14, 5).Value = 125
    $ws.Cells.Item(14, 6).Value = 103
    $ws.Cells.Item(15, 1).Value = 39
    $ws.Cells.Item(15, 2).Value = 102
    $ws.Cells.Item(15, 3).Value = 149
    $ws.Cells.Item(15, 4).Value = 83
    $ws.Cells.Item(15, 5).Value = 145
    $ws.Cells.Item(15, 6).Value = 134
    $ws.Cells.Item(16, 1).Value = 192
    $ws.Cells.Item(16, 2).Value = 62
    $ws.Cells.Item(16, 3).Value = 68
    $ws.Cells.Item(16, 4).Value = 77
    $ws.Cells.Item(16, 5).Value = 6
    $ws.Cells.Item(16, 6).Value = 124
    $ws.Cells.Item(17, 1).Value = 114
    $ws.Cells.Item(17, 2).Value = 73
    $ws.Cells.Item(17, 3).Value = 84
    $ws.Cells.Item(17, 4).Value = 92
    $ws.Cells.Item(17, 5).Value = 62
    $ws.Cells.Item(17, 6).Value = 45
    $ws.Cells.Item(18, 1).Value = 24
    $ws.Cells.Item(18, 2).Value = 27
    $ws.Cells.Item(18, 3).Value = 188
    $ws.Cells.Item(18, 4).Value = 98
    $ws.Cells.Item(18, 5).Value = 147
    $ws.Cells.Item(18, 6).Value = 61
    $ws.Cells.Item(19, 1).Value = 2
    $ws.Cells.Item(19, 2).Value = 123
    $ws.Cells.Item(19, 3).Value = 71
    $ws.Cells.Item(19, 4).Value = 128
    $ws.Cells.Item(19, 5).Value = 166
    $ws.Cells.Item(19, 6).Value = 20
    $ws.Cells.Item(20, 1).Value = 118
    $ws.Cells.Item(20, 2).Value = 116
    $ws.Cells.Item(20, 3).Value = 104
    $ws.Cells.Item(20, 4).Value = 90
    $ws.Cells.Item(20, 5).Value = 190
    $ws.Cells.Item(20, 6).Value = 36
    $ws.Cells.Item(21, 1).Value = 112
    $ws.Cells.Item(21, 2).Value = 125
    $ws.Cells.Item(21, 3).Value = 6
    $ws.Cells.Item(21, 4).Value = 127
    $ws.Cells.Item(21, 5).Value = 165
    $ws.Cells.Item(21, 6).Value = 78
    $ws.Cells.Item(22, 1).Value = 88
    $ws.Cells.Item(22, 2).Value = 17
    $ws.Cells.Item(22, 3).Value = 116
    $ws.Cells.Item(22, 4).Value = 172
    $ws.Cells.Item(22, 5).Value = 41
    $ws.Cells.Item(22, 6).Value = 19
    $ws.Cells.Item(23, 1).Value = 38
    $ws.Cells.Item(23, 2).Value = 133
    $ws.Cells.Item(23, 3).Value = 146
    $ws.Cells.Item(23, 4).Value = 125
    $ws.Cells.Item(23, 5).Value = 92
    $ws.Cells.Item(23, 6).Value = 62
    $ws.Cells.Item(24, 1).Value = 129
    $ws.Cells.Item(24, 2).Value = 112
    $ws.Cells.Item(24, 3).Value = 132
    $ws.Cells.Item(24, 4).Value = 25
    $ws.Cells.Item(24, 5).Value = 133
    $ws.Cells.Item(24, 6).Value = 196
    $ws.Cells.Item(25, 1).Value = 187
    $ws.Cells.Item(25, 2).Value = 165
    $ws.Cells.Item(25, 3).Value = 36
    $ws.Cells.Item(25, 4).Value = 188
    $ws.Cells.Item(25, 5).Value = 40
    $ws.Cells.Item(25, 6).Value = 43
    $ws.Cells.Item(26, 1).Value = 2
    $ws.Cells.Item(26, 2).Value = 60
    $ws.Cells.Item(26, 3).Value = 55
    $ws.Cells.Item(26, 4).Value = 116
    $ws.Cells.Item(26, 5).Value = 78
    $ws.Cells.Item(26, 6).Value = 37

# --- New Error1..Error5 / ErrorAverage columns start out zeroed, like the other computed columns ---
$ws.Range("V2:AA26").Value = 0
